# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values for rows 4 and 5
# on the zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-13 00:17:49"
$wsZhCn.Range("H4").Value = "2016-03-13 00:18:13"
$wsZhCn.Range("E5").Value = "2016-03-13 00:17:49"
$wsZhCn.Range("H5").Value = "2016-03-13 00:18:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-13 00:17:53"
$wsDeDe.Range("H4").Value = "2016-03-13 00:18:19"
$wsDeDe.Range("E5").Value = "2016-03-13 00:17:53"
$wsDeDe.Range("H5").Value = "2016-03-13 00:18:19"
